$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-16 03:49:09"
$ws.Range("E3").Value = "2026-02-16 03:49:11"
$ws.Range("E4").Value = "2026-02-16 03:49:14"
$ws.Range("J4").Value = "1014.7 hPa"
$ws.Range("O4").Value = "10.5 °C"
$ws.Range("E5").Value = "2026-02-16 03:49:17"
$ws.Range("I5").Value = "1.6 mm"
$ws.Range("E6").Value = "2026-02-16 03:49:19"
$ws.Range("J6").Value = "1014.8 hPa"
$ws.Range("L6").Value = "9.4 km/h - 302º 3:22 TU"
$ws.Range("O6").Value = "6.9 °C"
$ws.Range("E7").Value = "2026-02-16 03:49:22"
$ws.Range("J7").Value = "1015.1 hPa"
$ws.Range("M7").Value = "13.4 °C 3:26 TU"
$ws.Range("O7").Value = "13.1 °C"
$ws.Range("E8").Value = "2026-02-16 03:49:25"
$ws.Range("J8").Value = "1015.0 hPa"
$ws.Range("E9").Value = "2026-02-16 03:49:28"
$ws.Range("N9").Value = "4.2 °C 3:28 TU"
$ws.Range("O9").Value = "5.5 °C"
$ws.Range("E10").Value = "2026-02-16 03:49:30"
$ws.Range("M10").Value = "4.5 °C 3:29 TU"
$ws.Range("O10").Value = "3.9 °C"
$ws.Range("E11").Value = "2026-02-16 03:49:33"
$ws.Range("E12").Value = "2026-02-16 03:49:36"
$ws.Range("E13").Value = "2026-02-16 03:49:38"
$ws.Range("O13").Value = "1.5 °C"
$ws.Range("E14").Value = "2026-02-16 03:49:41"
$ws.Range("E15").Value = "2026-02-16 03:49:44"
$ws.Range("N15").Value = "4.6 °C 3:28 TU"
$ws.Range("O15").Value = "5.7 °C"
$ws.Range("E16").Value = "2026-02-16 03:49:47"
$ws.Range("H16").Value = "79%"
$ws.Range("I16").Value = "0.9 mm"
$ws.Range("M16").Value = "-0.3 °C 3:29 TU"
$ws.Range("O16").Value = "-0.9 °C"
$ws.Range("E17").Value = "2026-02-16 03:49:49"
$ws.Range("E18").Value = "2026-02-16 03:49:52"
$ws.Range("H18").Value = "98%"
$ws.Range("J18").Value = "1015.2 hPa"
$ws.Range("E19").Value = "2026-02-16 03:49:55"
$ws.Range("E20").Value = "2026-02-16 03:49:58"
$ws.Range("H20").Value = "89%"
$ws.Range("N20").Value = "-1.5 °C 3:09 TU"
$ws.Range("O20").Value = "-1.0 °C"
$ws.Range("E21").Value = "2026-02-16 03:50:00"
$ws.Range("H21").Value = "80%"
$ws.Range("J21").Value = "1017.0 hPa"
$ws.Range("E22").Value = "2026-02-16 03:50:03"
$ws.Range("I22").Value = "0.5 mm"
$ws.Range("N22").Value = "-6.5 °C 3:22 TU"
$ws.Range("E23").Value = "2026-02-16 03:50:06"
$ws.Range("E24").Value = "2026-02-16 03:50:09"
$ws.Range("H24").Value = "70%"
$ws.Range("J24").Value = "1018.3 hPa"
$ws.Range("O24").Value = "10.6 °C"
$ws.Range("E25").Value = "2026-02-16 03:50:11"
$ws.Range("E26").Value = "2026-02-16 03:50:14"
$ws.Range("E27").Value = "2026-02-16 03:50:17"
$ws.Range("E28").Value = "2026-02-16 03:50:19"
$ws.Range("L28").Value = "9.4 km/h - 210º 3:02 TU"
$ws.Range("E29").Value = "2026-02-16 03:50:22"
$ws.Range("E30").Value = "2026-02-16 03:50:25"
$ws.Range("H30").Value = "88%"
$ws.Range("J30").Value = "1014.8 hPa"
$ws.Range("E31").Value = "2026-02-16 03:50:27"
$ws.Range("H31").Value = "55%"
$ws.Range("J31").Value = "1013.3 hPa"
$ws.Range("O31").Value = "14.2 °C"
$ws.Range("E32").Value = "2026-02-16 03:50:30"
$ws.Range("H32").Value = "80%"
$ws.Range("E33").Value = "2026-02-16 03:50:33"
$ws.Range("H33").Value = "70%"
$ws.Range("N33").Value = "3.4 °C 3:29 TU"
$ws.Range("O33").Value = "5.1 °C"
$ws.Range("E34").Value = "2026-02-16 03:50:36"
$ws.Range("H34").Value = "65%"
$ws.Range("O34").Value = "3.2 °C"
$ws.Range("E35").Value = "2026-02-16 03:50:39"
$ws.Range("J35").Value = "1019.4 hPa"
$ws.Range("E36").Value = "2026-02-16 03:50:42"
$ws.Range("H36").Value = "89%"
$ws.Range("J36").Value = "1014.7 hPa"
$ws.Range("L36").Value = "9.4 km/h - 75º 3:20 TU"
$ws.Range("O36").Value = "7.1 °C"
$ws.Range("E37").Value = "2026-02-16 03:50:44"
$ws.Range("N37").Value = "1.4 °C 3:29 TU"
$ws.Range("O37").Value = "2.1 °C"
$ws.Range("E38").Value = "2026-02-16 03:50:47"
$ws.Range("H38").Value = "93%"
$ws.Range("O38").Value = "5.6 °C"
$ws.Range("E39").Value = "2026-02-16 03:50:50"
$ws.Range("H39").Value = "74%"
$ws.Range("E40").Value = "2026-02-16 03:50:53"
$ws.Range("H40").Value = "94%"
$ws.Range("N40").Value = "2.3 °C 3:02 TU"
$ws.Range("O40").Value = "3.3 °C"
$ws.Range("E41").Value = "2026-02-16 03:50:55"
$ws.Range("E42").Value = "2026-02-16 03:50:58"
$ws.Range("H42").Value = "95%"
$ws.Range("E43").Value = "2026-02-16 03:51:01"
$ws.Range("H43").Value = "99%"
$ws.Range("O43").Value = "3.6 °C"
$ws.Range("E44").Value = "2026-02-16 03:51:04"
$ws.Range("H44").Value = "87%"
$ws.Range("M44").Value = "0.3 °C 3:29 TU"
$ws.Range("O44").Value = "-0.4 °C"
$ws.Range("E45").Value = "2026-02-16 03:51:06"
$ws.Range("J45").Value = "1020.0 hPa"
$ws.Range("N45").Value = "2.9 °C 3:25 TU"
$ws.Range("E46").Value = "2026-02-16 03:51:09"
$ws.Range("J46").Value = "1018.8 hPa"
$ws.Range("O46").Value = "12.5 °C"
$ws.Range("G23").Value = "209 cm"
$ws.Range("H23").Value = "84%"
$ws.Range("I23").Value = "0.6 mm"
$ws.Range("K23").Value = "0.0 MJ/m2"
$ws.Range("L23").Value = "49.3 km/h - 318º 3:21 TU"
$ws.Range("M23").Value = "-0.7 °C 1:40 TU"
$ws.Range("N23").Value = "-1.5 °C 0:16 TU"
$ws.Range("O23").Value = "-1.1 °C"
